$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "0.0443") are preserved verbatim instead of being
# auto-coerced into numbers by Excel's usual Value-assignment heuristics.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.465.75'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.189.15'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +5.66%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.95'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.03'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +8.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.178.23'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +5.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.505'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.97%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +10.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.159'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.97%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.05'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +7.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000228'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.705.69'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.579.55'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.47%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.194.51'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.69%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.114'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '529.12'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +11.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.02'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +8.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.37'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.732'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +8.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.65'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +9.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.22'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +7.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.28'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.24'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +21.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.89'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +8.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.22'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +8.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.10'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.68'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.24%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +5.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '555.00'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.54'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.25'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +8.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.54'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0443'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +7.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0844'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +8.02%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +7.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.197.50'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +10.17%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.50'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +5.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.276'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +16.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.29'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +12.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.28'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +7.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.59'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0540'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.84%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.15%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +8.67%  '

Write-Output "Updated cryptos list"